# Update "想去人数" (want-to-go count) figures in column F across the
# three affected sheets to match the newly regenerated gh-pages output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 540
$ws1.Range("F4").Value = 14000
$ws1.Range("F5").Value = 241
$ws1.Range("F6").Value = 1807
$ws1.Range("F9").Value = 100
$ws1.Range("F11").Value = 559
$ws1.Range("F15").Value = 14164
$ws1.Range("F16").Value = 378
$ws1.Range("F17").Value = 637
$ws1.Range("F18").Value = 15034
$ws1.Range("F20").Value = 8388
$ws1.Range("F21").Value = 288
$ws1.Range("F23").Value = 38
$ws1.Range("F24").Value = 160
$ws1.Range("F30").Value = 39
$ws1.Range("F32").Value = 28
$ws1.Range("F36").Value = 12
$ws1.Range("F41").Value = 102
$ws1.Range("F42").Value = 5160

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 53

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 540
$ws4.Range("F4").Value = 14001
$ws4.Range("F5").Value = 241
$ws4.Range("F6").Value = 1807
$ws4.Range("F9").Value = 100
$ws4.Range("F11").Value = 559
$ws4.Range("F15").Value = 14164
$ws4.Range("F16").Value = 378
$ws4.Range("F17").Value = 637
$ws4.Range("F18").Value = 15034
$ws4.Range("F20").Value = 8388
$ws4.Range("F21").Value = 288
$ws4.Range("F23").Value = 38
$ws4.Range("F24").Value = 160
$ws4.Range("F30").Value = 39
$ws4.Range("F32").Value = 28
$ws4.Range("F35").Value = 53
$ws4.Range("F38").Value = 12
$ws4.Range("F43").Value = 102
$ws4.Range("F44").Value = 5160
